$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B1").Value = "business_type"
$f = $ws.Range("B1").Font
$f.Bold = $false
$f.Size = 16
$f.Name = "Arial"
$f.Color = 3880317
$ws.Rows("1:1").RowHeight = 20
$ws.Range("B1").Select()
Write-Output "done"
